$wb = $excel.ActiveWorkbook

# --- Update the Transactions sheet data (EPP Variable Installments T1 scenarios) ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()

$wsTrans.Range("A2").Value = 1454
$wsTrans.Range("E2").Value = 7.96
$wsTrans.Range("I2").Value = 7.96
$wsTrans.Range("A3").Value = 1453
$wsTrans.Range("A4").Value = 1449
$wsTrans.Range("A5").Value = 1448

# --- Repayment schedule sheet: selection moved to column Q ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
[void]$wsSchedule.Columns("Q").Select()

# --- Transactions sheet: selection moved to I5 ---
[void]$wsTrans.Range("I5").Select()

# --- NewLoanInput becomes the active tab ---
$wsInput = $wb.Worksheets.Item("NewLoanInput")
$wsInput.Activate()
